# Update countries & provincias Spain
# - Brasil overtakes Rusia in total cases -> swap their table rows (5 and 6)
# - Groenlandia / Seychelles swap order (tied case counts, list re-sorted)
# - Bonaire, San Eustaquio y Saba moves ahead of Sahara Occidental / San Bartolome
# - Refresh "Datos actualizados" timestamp
# - Refresh USA (row 4) totals

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 4: Estados Unidos - updated totals ---
$ws.Range("B4").Value = 1643304
$ws.Range("C4").Value = 22407
$ws.Range("D4").Value = 396417
$ws.Range("E4").Value = 1149325
$ws.Range("F4").Value = 0
$ws.Range("G4").Value = 1208
$ws.Range("H4").Value = 97562

# --- Row 5: becomes Brasil with its newly updated totals (overtakes Rusia) ---
$ws.Range("A5").Value = "Brasil"
$ws.Range("B5").Value = 330890
$ws.Range("C5").Value = 19969
$ws.Range("D5").Value = 125960
$ws.Range("E5").Value = 183882
$ws.Range("F5").Value = 0
$ws.Range("G5").Value = 966
$ws.Range("H5").Value = 21048

# --- Row 6: becomes Rusia, keeping its previous (unchanged) totals ---
$ws.Range("A6").Value = "Rusia"
$ws.Range("B6").Value = 326448
$ws.Range("C6").Value = 8894
$ws.Range("D6").Value = 99825
$ws.Range("E6").Value = 223374
$ws.Range("F6").Value = 0
$ws.Range("G6").Value = 150
$ws.Range("H6").Value = 3249

# --- Rows 209 / 210: Seychelles <-> Groenlandia (identical stats, order swap only) ---
$ws.Range("A209").Value = "Groenlandia"
$ws.Range("A210").Value = "Seychelles"

# --- Rows 214 / 215 / 216: Sahara Occidental, San Bartolome, Bonaire rotate ---
$ws.Range("A214").Value = "Bonaire, San Eustaquio y Saba"
$ws.Range("A215").Value = "Sahara Occidental"
$ws.Range("A216").Value = "San Bartolome"

# --- Timestamp refresh ---
$ws.Range("A1").Value = "Datos actualizados a 23 de Mayo de 2020 a las 00:35"
